$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.999.55"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.559.96"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'207.23"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").Value = "'22.14"
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("D11").Value = "'0.0858"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "1.781.49"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").Value = "1.541.71"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").Value = "'0.522"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "'62.01"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "26.981.26"
$ws.Range("D18").Value = "0.0₃0708"
$ws.Range("E18").Value = "  +2.82%  "
$ws.Range("D19").Value = "'217.39"
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("D20").Value = "'7.36"
$ws.Range("E20").Value = "  +1.70%  "
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").Value = "'4.10"
$ws.Range("E22").Value = "  +1.43%  "
$ws.Range("D23").Value = "'9.24"
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("D24").Value = "'1.93"
$ws.Range("E24").Value = "  -3.27%  "
$ws.Range("D25").Value = "'153.08"
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("D26").Value = "'6.65"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "'15.04"
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("D28").Value = "'0.105"
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'3.11"
$ws.Range("E33").Value = "  +3.70%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.422.47"
$ws.Range("E34").Value = "  +0.43%  "
$ws.Range("D35").Value = "'1.61"
$ws.Range("E35").Value = "  +2.53%  "
$ws.Range("D36").Value = "'1.03"
$ws.Range("E36").Value = "  +7.69%  "
$ws.Range("E37").Value = "  +0.96%  "
$ws.Range("D38").Value = "'0.0166"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("D39").Value = "'0.532"
$ws.Range("E39").Value = "  +1.74%  "
$ws.Range("D40").Value = "'0.809"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'1.01"
$ws.Range("E42").Value = "  +2.55%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.66"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "'2.33"
$ws.Range("E44").Value = "  +2.25%  "
$ws.Range("D45").Value = "'64.88"
$ws.Range("E45").Value = "  +1.95%  "
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("D47").Value = "1.695.27"
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").Value = "'87.47"
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("E51").Value = "  -0.10%  "
